$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Test 1"

# ---- Clear previous content/format ----
$ws.Cells.Clear()

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 18.42578125
$ws.Columns.Item(2).ColumnWidth = 37.7109375

# ---- Header row ----
$ws.Range("A1").Value = "Step Name"
$ws.Range("B1").Value = "Description"

$headerRange = $ws.Range("A1:K1")
$headerRange.Interior.Pattern = -4142
$headerRange.Interior.Color = 5855577
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 16
$headerRange.Font.Color = 16777215
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

# ---- Data rows ----
$steps = @(
  @("Step 1", "Navigate to URL."),
  @("Step 2", "Enter valid user name."),
  @("Step 3", "Enter valid password."),
  @("Step 4", "Click Login button."),
  @("Step 5", "Verify user navigated to Home Page."),
  @("Step 6", "Click SignOut button.")
)

for ($i = 0; $i -lt $steps.Length; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $steps[$i][0]
  $ws.Cells.Item($r, 2).Value = $steps[$i][1]
}

$dataRangeA = $ws.Range("A2:A7")
$dataRangeA.Font.Name = "Calibri"
$dataRangeA.Font.Size = 12
$dataRangeA.Font.Color = 0
$dataRangeA.HorizontalAlignment = -4108
$dataRangeA.VerticalAlignment = -4108

$dataRangeB = $ws.Range("B2:B7")
$dataRangeB.Font.Name = "Calibri"
$dataRangeB.Font.Size = 12
$dataRangeB.Font.Color = 0
$dataRangeB.HorizontalAlignment = -4131
$dataRangeB.VerticalAlignment = -4108

# ---- Row heights ----
for ($r = 1; $r -le 27; $r++) {
  $ws.Rows.Item($r).RowHeight = 20.1
}

# ---- Selection ----
$ws.Range("B18").Select()

Write-Host "done"
